$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.048.56"
$ws.Range("E2").Value = "  -0.49%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.272.76"
$ws.Range("E3").Value = "  -0.46%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - was Solana, now BNB
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "581.28"
$ws.Range("E5").Value = "  -1.06%  "

# Row 6 - was BNB, now Solana
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "185.91"
$ws.Range("E6").Value = "  -0.04%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.13%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.272.23"
$ws.Range("E9").Value = "  -0.45%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.10%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -2.09%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -1.96%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.841.43"
$ws.Range("E13").Value = "  -0.47%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.00%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "27.56"
$ws.Range("E15").Value = "  -4.74%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "68.040.83"
$ws.Range("E16").Value = "  -0.55%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -2.19%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.282.19"
$ws.Range("E18").Value = "  -0.41%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -1.94%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "13.60"
$ws.Range("E20").Value = "  -0.68%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "394.43"
$ws.Range("E21").Value = "  +2.62%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "7.63"
$ws.Range("E22").Value = "  -2.08%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "71.62"
$ws.Range("E23").Value = "  +0.21%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -4.02%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  -3.39%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "9.53"
$ws.Range("E28").Value = "  -2.75%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.52%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -2.63%  "

# Row 31 - was NEARProtocol, now EthereumClassic
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "22.71"
$ws.Range("E31").Value = "  -1.46%  "

# Row 32 - was EthereumClassic, now NEARProtocol
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "5.53"
$ws.Range("E32").Value = "  -5.75%  "

# Row 33 - Aptos
$ws.Range("D33").Value = "7.00"
$ws.Range("E33").Value = "  -3.31%  "

# Row 34 - Fetch.AI
$ws.Range("E34").Value = "  -4.67%  "

# Row 36 - Monero
$ws.Range("D36").Value = "163.91"
$ws.Range("E36").Value = "  +0.02%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "1.48"
$ws.Range("E37").Value = "  -4.48%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  +1.19%  "

# Row 39 - EnergySwap
$ws.Range("D39").Value = "26.81"
$ws.Range("E39").Value = "  +0.17%  "

# Row 40 - Mantle
$ws.Range("D40").Value = "0.812"
$ws.Range("E40").Value = "  -3.00%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  -1.71%  "

# Row 42 - RenderToken
$ws.Range("D42").Value = "6.47"
$ws.Range("E42").Value = "  -4.41%  "

# Row 43 - Hedera
$ws.Range("E43").Value = "  -0.73%  "

# Row 44 - Maker
$ws.Range("D44").Value = "2.632.93"
$ws.Range("E44").Value = "  -0.22%  "

# Row 45 - OKB
$ws.Range("D45").Value = "40.74"
$ws.Range("E45").Value = "  -1.34%  "

# Row 46 - dogwifhat
$ws.Range("E46").Value = "  -8.47%  "

# Row 47 - InjectiveProtocol
$ws.Range("D47").Value = "24.97"
$ws.Range("E47").Value = "  -3.13%  "

# Row 48 - Bittensor
$ws.Range("D48").Value = "333.34"
$ws.Range("E48").Value = "  -2.06%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  -2.72%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  +1.18%  "

# Row 51 - Stellar
$ws.Range("D51").Value = "0.102"
$ws.Range("E51").Value = "  -0.82%  "
